$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix typo: "Predicton" -> "Prediction"
$ws.Range("D2").Value = "Deep Learning Based Weather Prediction"

# Capitalize "The" and join the two line-broken sentences into one (removing the mid-paragraph break)
$ws.Range("K2").Value = "The conventional theory-driven numerical weather prediction (NWP) methods face many challenges, such`nas incomplete understanding of physical mechanisms, difficulties in obtaining useful knowledge from the deluge of observation data.Deep learning-based weather prediction (DLWP) is expected to be a strong supplement to the conventional method."

# Simplify H3 rich text (two runs: "Huawei Atlas AI" + " infrastructure") to a single plain-text run with the same visible content
$ws.Range("H3").Value = "Huawei Atlas AI infrastructure"

# Replace leading non-breaking space with a regular space and capitalize "The"
$ws.Range("F6").Value = " The continuous weather data of a particular region to predict the future weather conditions for the data analysis to predict the further weather conditions."

# Selection change
$ws.Range("H5").Select()

# Column width adjustments
$ws.Columns.Item(7).ColumnWidth = 18.140625
$ws.Columns.Item(8).ColumnWidth = 17.140625

# Row height adjustments
$ws.Rows.Item(3).RowHeight = 140.25
$ws.Rows.Item(4).RowHeight = 147
$ws.Rows.Item(6).RowHeight = 92.25
